$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6: update K6
$ws.Range("K6").Value = 1.388747888886706

# Row 7: update J7, add K7
$ws.Range("J7").Value = 0.444773652920949
$ws.Range("K7").Value = 0.2348700177716323

# Row 8: update I8, add J8
$ws.Range("I8").Value = 0.4487415504340581
$ws.Range("J8").Value = 0.2388379152847414

# Row 9: update H9, add I9
$ws.Range("H9").Value = 0.5843816406042994
$ws.Range("I9").Value = 0.3744780054549828

# Row 10: update G10, add H10
$ws.Range("G10").Value = 0.3435754587486348
$ws.Range("H10").Value = 0.1336718235993181

# Row 11: update F11, add G11
$ws.Range("F11").Value = 0.2982442434965384
$ws.Range("G11").Value = 0.08834060834722172

# Row 12: update E12, add F12
$ws.Range("E12").Value = 0.2313828215604846
$ws.Range("F12").Value = 0.02147918641116785

# Row 13: update D13, add E13
$ws.Range("D13").Value = 0.201796619203768
$ws.Range("E13").Value = -0.00810701594554874

# Row 14: update C14, add D14
$ws.Range("C14").Value = 0.1836459624741271
$ws.Range("D14").Value = -0.02625767267518964

# Row 15: update B15, add C15
$ws.Range("B15").Value = 0.1656141382254278
$ws.Range("C15").Value = -0.04428949692388896

# Row 16: add B16
$ws.Range("B16").Value = -0.09587373626955231
